# Auto-generated Excel COM-interop script
# Applies the numeric cell updates described in the commit diff
# to the Odin_Profits workbook (sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ALC_updates = @{
    "H28" = 6672.952
    "I28" = 993.625
    "J28" = 10167.923
    "K28" = 993.625
    "L28" = 10167.923
    "M28" = -508.625
    "N28" = -11137.923
    "H42" = 739.8570999999999
    "I42" = 846.5
    "J42" = 100
    "K42" = 2539.5
    "L42" = 300
    "M42" = -2309.5
    "N42" = -760
    "H43" = 9471
    "I43" = 8871.75
    "K43" = 8871.75
    "M43" = -8802.75
    "H58" = 2304.8333
    "I58" = 95
    "J58" = 6724.5
    "K58" = 285
    "L58" = 20173.5
    "M58" = -135
    "N58" = -20473.5
    "H61" = 682.8
    "I61" = 403.5
    "J61" = 1800
    "K61" = 1210.5
    "L61" = 5400
    "M61" = -1038.5
    "N61" = -5744
    "H76" = 35720244
    "I76" = 50006810
    "K76" = 50006810
    "M76" = -50006495
    "H79" = 35720244
    "I79" = 50006810
    "K79" = 50006810
    "M79" = -50005718
    "H88" = 2024.8096
    "I88" = 985.1667
    "J88" = 3411
    "K88" = 985.1667
    "L88" = 3411
    "M88" = -579.1667
    "N88" = -4223
    "H91" = 2024.8096
    "I91" = 985.1667
    "J91" = 3411
    "K91" = 985.1667
    "L91" = 3411
    "M91" = 418.8333
    "N91" = -6219
    "H106" = 5022.625
    "I106" = 4152.3335
    "K106" = 4152.3335
    "M106" = -3521.3335
    "H113" = 3203.889
    "I113" = 3430.8333
    "K113" = 3430.8333
    "M113" = -176.8332999999998
    "H121" = 0
    "J121" = 0
    "L121" = 0
    "H132" = 351951.84
    "I132" = 401750.75
    "K132" = 1205252.25
    "M132" = -1202722.25
    "H135" = 5179
    "J135" = 11664
    "L135" = 104976
    "N135" = -110046
}
foreach ($cellRef in $ALC_updates.Keys) {
    $ws.Range($cellRef).Value = $ALC_updates[$cellRef]
}

$ws.Range("N121").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ARM_updates = @{
    "H32" = 9859281
    "I32" = 8929721
    "K32" = 8929721
    "M32" = -8929434
    "H61" = 4086.024
    "I61" = 3568.3872
    "J61" = 5544.8184
    "K61" = 3568.3872
    "L61" = 5544.8184
    "M61" = -3356.3872
    "N61" = -5968.8184
    "H63" = 3095.7693
    "J63" = 1499
    "L63" = 1499
    "N63" = -2871
    "H66" = 3095.7693
    "J66" = 1499
    "L66" = 7495
    "N66" = -14359
    "H122" = 2872.25
    "I122" = 2115.6428
    "K122" = 6346.928400000001
    "M122" = -3896.928400000001
    "H132" = 3343555.8
    "I132" = 8073890
    "K132" = 24221670
    "M132" = -24219140
    "H136" = 4086.024
    "I136" = 3568.3872
    "J136" = 5544.8184
    "K136" = 10705.1616
    "L136" = 16634.4552
    "M136" = -8155.161599999999
    "N136" = -21734.4552
}
foreach ($cellRef in $ARM_updates.Keys) {
    $ws.Range($cellRef).Value = $ARM_updates[$cellRef]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$CRP_updates = @{
    "H5" = 1561.25
    "I5" = 225
    "K5" = 225
    "M5" = -113
    "H31" = 6541.591
    "I31" = 1976.6666
    "K31" = 1976.6666
    "M31" = -1681.6666
    "H34" = 6541.591
    "I34" = 1976.6666
    "K34" = 1976.6666
    "M34" = -1774.6666
    "H99" = 15876540
    "I99" = 18521796
    "J99" = 5000
    "K99" = 18521796
    "L99" = 5000
    "M99" = -18520298
    "N99" = -7996
    "H126" = 15876540
    "I126" = 18521796
    "J126" = 5000
    "K126" = 55565388
    "L126" = 15000
    "M126" = -55562918
    "N126" = -19940
    "H132" = 11526.833
    "I132" = 8884.1
    "K132" = 26652.3
    "M132" = -24122.3
}
foreach ($cellRef in $CRP_updates.Keys) {
    $ws.Range($cellRef).Value = $CRP_updates[$cellRef]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$CUL_updates = @{
    "H23" = 207.25
    "I23" = 49
    "J23" = 260
    "K23" = 147
    "L23" = 780
    "M23" = 88
    "N23" = -1250
    "H97" = 479.4
    "J97" = 499
    "L97" = 1497
    "N97" = -2489
    "H113" = 1015.38464
    "I113" = 576
    "J113" = 1052
    "K113" = 1728
    "L113" = 3156
    "M113" = 442
    "N113" = -7496
    "H137" = 2162.84
    "I137" = 801.1
    "J137" = 3070.6667
    "K137" = 2403.3
    "L137" = 9212.000100000001
    "M137" = 2696.7
    "N137" = -19412.0001
}
foreach ($cellRef in $CUL_updates.Keys) {
    $ws.Range($cellRef).Value = $CUL_updates[$cellRef]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$GSM_updates = @{
    "H122" = 4211.6313
    "I122" = 1958.8462
    "J122" = 9092.666999999999
    "K122" = 5876.5386
    "L122" = 27278.001
    "M122" = -3426.5386
    "N122" = -32178.001
    "H126" = 30008742
    "I126" = 50003544
    "J126" = 10013940
    "K126" = 150010632
    "L126" = 30041820
    "M126" = -150008162
    "N126" = -30046760
    "H132" = 26320474
    "I132" = 37041772
    "K132" = 111125316
    "M132" = -111122786
}
foreach ($cellRef in $GSM_updates.Keys) {
    $ws.Range($cellRef).Value = $GSM_updates[$cellRef]
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$LTW_updates = @{
    "H16" = 3479.7896
    "I16" = 905.0769
    "K16" = 905.0769
    "M16" = -735.0769
    "H40" = 4529.3887
    "I40" = 4725.423
    "J40" = 4019.7
    "K40" = 4725.423
    "L40" = 4019.7
    "M40" = -4589.423
    "N40" = -4291.7
    "H100" = 5036.625
    "I100" = 5963.1665
    "J100" = 2257
    "K100" = 5963.1665
    "L100" = 2257
    "M100" = -5422.1665
    "N100" = -3339
    "H122" = 3669.2083
    "I122" = 3241.7693
    "J122" = 4174.364
    "K122" = 9725.3079
    "L122" = 12523.092
    "M122" = -7275.3079
    "N122" = -17423.092
}
foreach ($cellRef in $LTW_updates.Keys) {
    $ws.Range($cellRef).Value = $LTW_updates[$cellRef]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$WVR_updates = @{
    "H107" = 7408240
    "I107" = 11111674
    "J107" = 1370.6666
    "K107" = 33335022
    "L107" = 4111.9998
    "M107" = -33333102
    "H136" = 13901658
    "I136" = 19240448
    "K136" = 57721344
    "M136" = -57718794
}
foreach ($cellRef in $WVR_updates.Keys) {
    $ws.Range($cellRef).Value = $WVR_updates[$cellRef]
}
